$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fgf2"
$ws.Range("C2").Value = "Fgfr2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.005243333333333
$ws.Range("H2").Value = 3.01573
$ws.Range("I2").Value = 0.07224874268505826
$ws.Range("J2").Value = 0.07224874268505825
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.033285
$ws.Range("N2").Value = 0.099855
$ws.Range("O2").Value = 0.007684499559038781
$ws.Range("P2").Value = 0.007684499559038781
$ws.Range("Q2").Value = 0.03345952435
$ws.Range("R2").Value = 0.30113571915
$ws.Range("S2").Value = 0.0005551954313044366
$ws.Range("T2").Value = 0.0005551954313044365

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fgf2"
$ws.Range("C3").Value = "Fgfr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.005243333333333
$ws.Range("H3").Value = 3.01573
$ws.Range("I3").Value = 0.07224874268505826
$ws.Range("J3").Value = 0.07224874268505825
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.140873
$ws.Range("N3").Value = 12.422619
$ws.Range("O3").Value = 0.9560023056192156
$ws.Range("P3").Value = 0.9560023056192157
$ws.Range("Q3").Value = 4.16258497743
$ws.Range("R3").Value = 37.46326479687
$ws.Range("S3").Value = 0.06906996458500514
$ws.Range("T3").Value = 0.06906996458500514

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fgf2"
$ws.Range("C4").Value = "Fgfr2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.005243333333333
$ws.Range("H4").Value = 3.01573
$ws.Range("I4").Value = 0.07224874268505826
$ws.Range("J4").Value = 0.07224874268505825
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.1572886666666667
$ws.Range("N4").Value = 0.471866
$ws.Range("O4").Value = 0.03631319482174546
$ws.Range("P4").Value = 0.03631319482174546
$ws.Range("Q4").Value = 0.1581133835755555
$ws.Range("R4").Value = 1.42302045218
$ws.Range("S4").Value = 0.002623582668748678
$ws.Range("T4").Value = 0.002623582668748678

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Fgf2"
$ws.Range("C5").Value = "Fgfr2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 10.25983933333333
$ws.Range("H5").Value = 30.779518
$ws.Range("I5").Value = 0.7373940889775011
$ws.Range("J5").Value = 0.737394088977501
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.033285
$ws.Range("N5").Value = 0.099855
$ws.Range("O5").Value = 0.007684499559038781
$ws.Range("P5").Value = 0.007684499559038781
$ws.Range("Q5").Value = 0.34149875221
$ws.Range("R5").Value = 3.07348876989
$ws.Range("S5").Value = 0.005666504551585411
$ws.Range("T5").Value = 0.00566650455158541

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Fgf2"
$ws.Range("C6").Value = "Fgfr2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 10.25983933333333
$ws.Range("H6").Value = 30.779518
$ws.Range("I6").Value = 0.7373940889775011
$ws.Range("J6").Value = 0.737394088977501
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 4.140873
$ws.Range("N6").Value = 12.422619
$ws.Range("O6").Value = 0.9560023056192156
$ws.Range("P6").Value = 0.9560023056192157
$ws.Range("Q6").Value = 42.484691679738
$ws.Range("R6").Value = 382.362225117642
$ws.Range("S6").Value = 0.7049504492124721
$ws.Range("T6").Value = 0.7049504492124721

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Fgf2"
$ws.Range("C7").Value = "Fgfr2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 10.25983933333333
$ws.Range("H7").Value = 30.779518
$ws.Range("I7").Value = 0.7373940889775011
$ws.Range("J7").Value = 0.737394088977501
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.1572886666666667
$ws.Range("N7").Value = 0.471866
$ws.Range("O7").Value = 0.03631319482174546
$ws.Range("P7").Value = 0.03631319482174546
$ws.Range("Q7").Value = 1.613756448954222
$ws.Range("R7").Value = 14.523808040588
$ws.Range("S7").Value = 0.02677713521344351
$ws.Range("T7").Value = 0.0267771352134435

# Row 8
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Fgf2"
$ws.Range("C8").Value = "Fgfr2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.383875
$ws.Range("H8").Value = 1.151625
$ws.Range("I8").Value = 0.02758982345723265
$ws.Range("J8").Value = 0.02758982345723265
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.033285
$ws.Range("N8").Value = 0.099855
$ws.Range("O8").Value = 0.007684499559038781
$ws.Range("P8").Value = 0.007684499559038781
$ws.Range("Q8").Value = 0.012777279375
$ws.Range("R8").Value = 0.114995514375
$ws.Range("S8").Value = 0.0002120139861910621
$ws.Range("T8").Value = 0.0002120139861910621

# Row 9
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Fgf2"
$ws.Range("C9").Value = "Fgfr2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.383875
$ws.Range("H9").Value = 1.151625
$ws.Range("I9").Value = 0.02758982345723265
$ws.Range("J9").Value = 0.02758982345723265
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 4.140873
$ws.Range("N9").Value = 12.422619
$ws.Range("O9").Value = 0.9560023056192156
$ws.Range("P9").Value = 0.9560023056192157
$ws.Range("Q9").Value = 1.589577622875
$ws.Range("R9").Value = 14.306198605875
$ws.Range("S9").Value = 0.02637593483674154
$ws.Range("T9").Value = 0.02637593483674154

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Fgf2"
$ws.Range("C10").Value = "Fgfr2"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.383875
$ws.Range("H10").Value = 1.151625
$ws.Range("I10").Value = 0.02758982345723265
$ws.Range("J10").Value = 0.02758982345723265
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.1572886666666667
$ws.Range("N10").Value = 0.471866
$ws.Range("O10").Value = 0.03631319482174546
$ws.Range("P10").Value = 0.03631319482174546
$ws.Range("Q10").Value = 0.06037918691666667
$ws.Range("R10").Value = 0.5434126822500001
$ws.Range("S10").Value = 0.001001874634300052
$ws.Range("T10").Value = 0.001001874634300052

# Row 11
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Fgf2"
$ws.Range("C11").Value = "Fgfr2"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 2.264687
$ws.Range("H11").Value = 6.794061
$ws.Range("I11").Value = 0.162767344880208
$ws.Range("J11").Value = 0.162767344880208
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.033285
$ws.Range("N11").Value = 0.099855
$ws.Range("O11").Value = 0.007684499559038781
$ws.Range("P11").Value = 0.007684499559038781
$ws.Range("Q11").Value = 0.075380106795
$ws.Range("R11").Value = 0.6784209611550001
$ws.Range("S11").Value = 0.001250785589957872
$ws.Range("T11").Value = 0.001250785589957872

# Row 12
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Fgf2"
$ws.Range("C12").Value = "Fgfr2"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 2.264687
$ws.Range("H12").Value = 6.794061
$ws.Range("I12").Value = 0.162767344880208
$ws.Range("J12").Value = 0.162767344880208
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 4.140873
$ws.Range("N12").Value = 12.422619
$ws.Range("O12").Value = 0.9560023056192156
$ws.Range("P12").Value = 0.9560023056192157
$ws.Range("Q12").Value = 9.377781251750999
$ws.Range("R12").Value = 84.40003126575901
$ws.Range("S12").Value = 0.1556059569849969
$ws.Range("T12").Value = 0.1556059569849969

# Row 13
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Fgf2"
$ws.Range("C13").Value = "Fgfr2"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 2.264687
$ws.Range("H13").Value = 6.794061
$ws.Range("I13").Value = 0.162767344880208
$ws.Range("J13").Value = 0.162767344880208
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.1572886666666667
$ws.Range("N13").Value = 0.471866
$ws.Range("O13").Value = 0.03631319482174546
$ws.Range("P13").Value = 0.03631319482174546
$ws.Range("Q13").Value = 0.3562095986473333
$ws.Range("R13").Value = 3.205886387826
$ws.Range("S13").Value = 0.005910602305253227
$ws.Range("T13").Value = 0.005910602305253227
